$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated s_val data (regen sval data to filter save games)
$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
    3 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 0, 18.67282528286833)
    4 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 0, 14.05633640148523)
    5 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
    6 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E - IP
    $ws.Cells.Item($row, 6).Value = $vals[4]  # F - Win
    $ws.Cells.Item($row, 7).Value = $vals[5]  # G - sum
}
